$d = $word.ActiveDocument

# --- 1. Remove the "Meta description: ..." paragraph that used to follow
#        the Heading1 title at the top of the document. ---
$findRange = $d.Content.Duplicate
$found = $findRange.Find.Execute("Meta description", $true, $false, $false,
                                  $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $findRange.Expand(4) | Out-Null   # wdParagraph: grow to the whole paragraph (incl. mark)
    $findRange.Delete()
}

# --- 2. Before the final paragraph (the image-prompt paragraph), insert a
#        new bold paragraph with the page title, and replace the final
#        paragraph's italic text with the meta-description text. The
#        replacement is done via InsertXML using the same run shape
#        (leading empty run + formatted run) used elsewhere in this
#        document, so the resulting OOXML matches the document's own
#        conventions. ---
$count = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($count)
$lastRange = $last.Range
$targetRange = $d.Range($lastRange.Start, $lastRange.End - 1)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Crazy Colt for Free: A Simple Slot Game with High Payouts</w:t></w:r></w:p><w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Crazy Colt, a simple online slot game with Wild and Scatter symbols and 50 paylines. Play for free and win big!</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$targetRange.InsertXML($xml)
